# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Femacal de La Calera - Bruselas (repollito)"
# at row 56 (pushing the existing rows 56-63 down to 57-64).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 56; everything currently at/after row 56 shifts down one row
# (old 56 -> 57, ... old 63 -> 64), and the sheet dimension grows from R63 to R64.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly data point.
$ws.Range("A56").Value = 3
$ws.Range("B56").Value = "Femacal de La Calera"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value = 44769
$ws.Range("E56").Value = 5
$ws.Range("F56").Value = 100112035
$ws.Range("G56").Value = "Bruselas (repollito)"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 85
$ws.Range("K56").Value = 14000
$ws.Range("L56").Value = 15000
$ws.Range("M56").Value = 14471
$ws.Range("N56").Value = "`$/malla 15 kilos"
$ws.Range("O56").Value = "Provincia de Quillota"
$ws.Range("P56").Value = 965
$ws.Range("Q56").Value = 15
$ws.Range("R56").Value = "Hortaliza"
